# Fix typo in the "comparison_title" column of the metadata sheet:
# "Obstructed defecation symdrome" -> "Obstructed defecation syndrome"
# (cell G2 on Sheet1, referenced via the shared-strings table).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$cell = $ws.Range("G2")
if ($cell.Value2 -eq "Obstructed defecation symdrome") {
    $cell.Value = "Obstructed defecation syndrome"
}
